$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Rename "ROI рекламы, млрд" (A21) to "Выручка от рекламы, млрд"
$ws.Range("A21").Value = "Выручка от рекламы, млрд"

# 2. C21 formula: =C20*14.95 -> =C20*B5 (reference the Выручка в 2018 cell instead of the literal)
$ws.Range("C21").Formula = "=C20*B5"

# 3. C22 formula: =C21*(1+B13) -> =C21*(1+0.7*B13)
$ws.Range("C22").Formula = "=C21*(1+0.7*B13)"
$ws.Range("C22").ClearFormats()

# 4. C23 number format: 0.000 -> 0.0
$ws.Range("C23").NumberFormat = "0.0"

# 5. Remove old rows 24 ("Для всего MasterCard") and 25 ("Прирост к чистой прибыли, в млн.$"),
#    which shifts the trailing blank rows (28, 29) up to (26, 27).
$ws.Rows("24:25").Delete()
